# Insert a new weekly price record as row 14 (pushing the existing
# rows 14-59 down to 15-60), matching the author's commit
# "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 14..59 down by one row.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new data point.
$ws.Cells.Item(14, 1).Value  = 1
$ws.Cells.Item(14, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(14, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(14, 4).Value  = 44592
$ws.Cells.Item(14, 5).Value  = 15
$ws.Cells.Item(14, 6).Value  = 100114001
$ws.Cells.Item(14, 7).Value  = "Papa"
$ws.Cells.Item(14, 8).Value  = "Patagonia"
$ws.Cells.Item(14, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(14, 10).Value = 1000
$ws.Cells.Item(14, 11).Value = 10000
$ws.Cells.Item(14, 12).Value = 11000
$ws.Cells.Item(14, 13).Value = 10500
$ws.Cells.Item(14, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(14, 15).Value = "Región del Maule"
$ws.Cells.Item(14, 16).Value = 420
$ws.Cells.Item(14, 17).Value = 25
$ws.Cells.Item(14, 18).Value = "Hortaliza"

# Give the new date cell the same date-time number format as the rest
# of column D.
$ws.Cells.Item(14, 4).NumberFormat = $ws.Cells.Item(15, 4).NumberFormat
